# Courier invoice / misc-charges sheet changes #CRM-1203
# The "Product" column (column D) is not required in the misc charges
# sheet - everything billed to the partner was listed as "Service" in
# that column, so it added no value. Remove the entire column; this
# shifts the old "Description" (E) and "Charge" (F) columns left into
# D and E respectively, and Excel automatically fixes up the merged
# cells (B2:F2 -> B2:E2, B12:D13 -> B12:C13) and column widths.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(4).Delete()
